# Add "2022-Q1" sheet data, reusing the existing "总计" worksheet in place
# (it keeps its physical sheet / sheetId / rId), then create a brand new
# "总计" worksheet after it with the updated totals table.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$text) {
    # Force a value to be stored as text even when it "looks like" a number,
    # without leaving a residual style on the cell itself.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

function Copy-HeaderStyle($srcRange, $dstRange) {
    # Clone the existing bordered/bold/centered header style (already used
    # by this workbook) onto a new cell without creating a new style entry.
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)   # xlPasteFormats
}

# ---------------------------------------------------------------------
# Step 1: rename the current "总计" sheet to "2022-Q1" and give it the
# fund-holdings table layout (like the other quarterly sheets).
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# New header cells E1:H1 need the same style already used by B1:D1.
Copy-HeaderStyle $q1.Range("B1") $q1.Range("E1")
Copy-HeaderStyle $q1.Range("B1") $q1.Range("F1")
Copy-HeaderStyle $q1.Range("B1") $q1.Range("G1")
Copy-HeaderStyle $q1.Range("B1") $q1.Range("H1")

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row 2 holds the single fund entry for 2022-Q1.
Set-TextValue $q1.Range("B2") "004397"
Set-TextValue $q1.Range("C2") "长盛信息安全量化策略灵活配置混合"
Set-TextValue $q1.Range("D2") "4.21"
Set-TextValue $q1.Range("E2") "29.75"
Set-TextValue $q1.Range("F2") "0.68"
Set-TextValue $q1.Range("G2") "0.0286"
$q1.Range("H2").Value = 9

# Drop the old total rows (3:5) that used to live on this sheet; this also
# shrinks the worksheet dimension back down to A1:H2.
$q1.Range("A3:D5").EntireRow.Delete()

# ---------------------------------------------------------------------
# Step 2: insert a brand-new "总计" worksheet right after "2022-Q1" and
# rebuild the summary table (now including the 2022-Q1 row).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

Copy-HeaderStyle $q1.Range("B1") $total.Range("B1")
Copy-HeaderStyle $q1.Range("B1") $total.Range("C1")
Copy-HeaderStyle $q1.Range("B1") $total.Range("D1")

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

Copy-HeaderStyle $q1.Range("A2") $total.Range("A2")
Copy-HeaderStyle $q1.Range("A2") $total.Range("A3")
Copy-HeaderStyle $q1.Range("A2") $total.Range("A4")
Copy-HeaderStyle $q1.Range("A2") $total.Range("A5")
Copy-HeaderStyle $q1.Range("A2") $total.Range("A6")

$total.Range("A2").Value = 0
Set-TextValue $total.Range("B2") "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.03

$total.Range("A3").Value = 1
Set-TextValue $total.Range("B3") "2021-Q4"
$total.Range("C3").Value = 4
$total.Range("D3").Value = 1.07

$total.Range("A4").Value = 2
Set-TextValue $total.Range("B4") "2021-Q3"
$total.Range("C4").Value = 3
$total.Range("D4").Value = 0.19

$total.Range("A5").Value = 3
Set-TextValue $total.Range("B5") "2021-Q2"
$total.Range("C5").Value = 6
$total.Range("D5").Value = 8.039999999999999

$total.Range("A6").Value = 4
Set-TextValue $total.Range("B6") "2020-Q4"
$total.Range("C6").Value = 3
$total.Range("D6").Value = 0.27

# Keep the originally-active sheet selected, like before the edit.
$wb.Worksheets.Item(1).Activate()
